# Textbox response formatting fix
# Renames sheets (regenerated task-order timestamps) and updates the
# stimulus-file / response-label values within each sheet.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order ids refreshed) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-1651168684318216"
$wb.Worksheets.Item(2).Name = "NB_TO-16511686878871334"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168687888048"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511686879350808"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511686879990811"

# --- Sheet 1: GNG_TO... ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651168684277213.csv"
$ws1.Range("B3").Value = "GNG_stims-16511686843002436.csv"
$ws1.Range("B4").Value = "go_stims-16511686843012083.csv"
$ws1.Range("B5").Value = "GNG_stims-1651168684317216.csv"

# --- Sheet 2: NB_TO... ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16511686865591183.csv"
$ws2.Range("B3").Value = "ZB-match_0-16511686843772082.csv"
$ws2.Range("B4").Value = "TB-1651168687230525.csv"
$ws2.Range("B5").Value = "OB-16511686847159374.csv"
$ws2.Range("B6").Value = "ZB-match_9-16511686845589345.csv"
$ws2.Range("B7").Value = "ZB-match_2-1651168684598935.csv"
$ws2.Range("B8").Value = "OB-16511686858579407.csv"
$ws2.Range("B9").Value = "OB-16511686854339457.csv"
$ws2.Range("B10").Value = "TB-16511686878650455.csv"

# --- Sheet 3: RS_TO... ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO... ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1651168687903045.csv"
$ws4.Range("B3").Value = "ZM_stims-16511686878900456.csv"
$ws4.Range("B4").Value = "MM_stims-16511686879190795.csv"
$ws4.Range("B5").Value = "ZM_stims-1651168687903045.csv"
$ws4.Range("B6").Value = "MM_stims-16511686879350808.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686879201021.csv"

# --- Sheet 5: vSAT_TO... ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511686879510455.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511686879830818.csv"
$ws5.Range("B4").Value = "vSAT_stims-165116868796708.csv"
$ws5.Range("B5").Value = "SAT_stims-1651168687938046.csv"
